$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Overview sheet: handback status text + column widths
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E1").ColumnWidth = 29.1
$wsOverview.Range("F1").ColumnWidth = 29.1

# ---------------------------------------------------------------------------
# 2) zh-cn sheet: target file / handback file / handback datetime
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$zhUrl = "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/a7ce1bdb2691afb04a5f6fa6273f66ca22aa703f/e2e/2f53bac2-edd1-4225-bb99-0580987f80b5.md"

# "Status" column shares the same underlying text as the Overview summary cells
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"

$wsZh.Hyperlinks.Add($wsZh.Range("J2"), $zhUrl, "", "", "2f53bac2-edd1-4225-bb99-0580987f80b5.md")
$wsZh.Range("J2").Style = "HyperLink"
$wsZh.Hyperlinks.Add($wsZh.Range("J3"), $zhUrl, "", "", "2f53bac2-edd1-4225-bb99-0580987f80b5.md")
$wsZh.Range("J3").Style = "HyperLink"

$wsZh.Range("K2").Value = "2f53bac2-edd1-4225-bb99-0580987f80b5.195f8633f6b48e9107b1c824c21c8e9b7e59076f.zh-cn.xlf"
$wsZh.Range("K3").Value = "2f53bac2-edd1-4225-bb99-0580987f80b5.195f8633f6b48e9107b1c824c21c8e9b7e59076f.zh-cn.xlf"

$wsZh.Range("L2").Value = "2017-02-22 07:50:57"
$wsZh.Range("L3").Value = "2017-02-22 07:50:57"

$wsZh.Range("C1").ColumnWidth = 29.1
$wsZh.Range("J1").ColumnWidth = 39.1
$wsZh.Range("K1").ColumnWidth = 39.1

# ---------------------------------------------------------------------------
# 3) de-de sheet: target file / handback file / handback datetime
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$deUrl = "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/a7ce1bdb2691afb04a5f6fa6273f66ca22aa703f/e2e/2f53bac2-edd1-4225-bb99-0580987f80b5.md"

$wsDe.Hyperlinks.Add($wsDe.Range("J2"), $deUrl, "", "", "2f53bac2-edd1-4225-bb99-0580987f80b5.md")
$wsDe.Range("J2").Style = "HyperLink"
$wsDe.Hyperlinks.Add($wsDe.Range("J3"), $deUrl, "", "", "2f53bac2-edd1-4225-bb99-0580987f80b5.md")
$wsDe.Range("J3").Style = "HyperLink"

$wsDe.Range("K2").Value = "2f53bac2-edd1-4225-bb99-0580987f80b5.195f8633f6b48e9107b1c824c21c8e9b7e59076f.de-de.xlf"
$wsDe.Range("K3").Value = "2f53bac2-edd1-4225-bb99-0580987f80b5.195f8633f6b48e9107b1c824c21c8e9b7e59076f.de-de.xlf"

$wsDe.Range("L2").Value = "2017-02-22 07:51:20"
$wsDe.Range("L3").Value = "2017-02-22 07:51:20"

$wsDe.Range("C1").ColumnWidth = 29.1
$wsDe.Range("J1").ColumnWidth = 39.1
$wsDe.Range("K1").ColumnWidth = 39.1

Write-Host "Handback report generated"
